$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows (13/05 - 27/05 updates, rows 256-269)
$data = @(
    @(256, 44330, 3, 11, 109.2353525322741),
    @(257, 44331, 4, 14, 139.0268123138034),
    @(258, 44332, 2, 14, 139.0268123138034),
    @(259, 44333, 2, 12, 119.1658391261172),
    @(260, 44334, 0, 11, 109.2353525322741),
    @(261, 44335, 0, 11, 109.2353525322741),
    @(262, 44336, 0, 11, 109.2353525322741),
    @(263, 44337, 1, 9, 89.37437934458789),
    @(264, 44338, 0, 5, 49.65243296921549),
    @(265, 44339, 0, 3, 29.7914597815293),
    @(266, 44340, 0, 1, 9.930486593843098),
    @(267, 44341, 0, 1, 9.930486593843098),
    @(268, 44342, 1, 2, 19.8609731876862),
    @(269, 44343, 0, 2, 19.8609731876862)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]

    # Copy formatting (style) from the row above, matching the existing pattern
    $prev = $r - 1
    $ws.Range("A$prev").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
}
